{"js": "// Apply the four textual edits described in the commit\n// \"Anpassungen Einleitung , Organigramm\".\n//\n// Each edit is performed as a targeted search-and-replace on a unique,\n// sufficiently long substring so that there is no risk of an accidental\n// match elsewhere in the document.\n\nconst doc = context.document;\nconst body = doc.body;\n\n// 1) \"Zuverl\u00e4ssigkeit\" table / Standards cell:\n//    \"Zu aufsetzende Dokumente ... abgegeben, der Zust\u00e4ndige stellt...\"\n//      -> \"Aufzusetzende Dokumente ... abgegeben. Der Zust\u00e4ndige stellt...\"\nconst search1 = body.search(\n  \"Zu aufsetzende Dokumente werden vom jeweiligen aufgetragenen PM zum vereinbarten Zeitpunkt abgegeben, der Zust\u00e4ndige stellt zudem sicher, dass\",\n  { matchCase: true }\n);\nsearch1.load(\"items\");\nawait context.sync();\nif (search1.items.length > 0) {\n  search1.items[0].insertText(\n    \"Aufzusetzende Dokumente werden vom jeweiligen aufgetragenen PM zum vereinbarten Zeitpunkt abgegeben. Der Zust\u00e4ndige stellt zudem sicher, dass\",\n    \"Replace\"\n  );\n  await context.sync();\n}\n\n// 2) \"Engagement\" table / Indikatoren cell:\n//    \"...fragen sie nach Hilfe und nehmen diese an.\"\n//      -> \"...fragen sie nach Hilfe und nehmen diese gegebenenfalls an.\"\nconst search2 = body.search(\n  \"fragen sie nach Hilfe und nehmen diese an.\",\n  { matchCase: true }\n);\nsearch2.load(\"items\");\nawait context.sync();\nif (search2.items.length > 0) {\n  search2.items[0].insertText(\n    \"fragen sie nach Hilfe und nehmen diese gegebenenfalls an.\",\n    \"Replace\"\n  );\n  await context.sync();\n}\n\n// 3) \"Engagement\" table / Standards cell, paragraph 1:\n//    \"...diskutiert und festgelegt, jedes PM macht sich Gedanken zu\n//    Sinnhaftigkeit und Effizienz.\"\n//      -> \"...diskutiert und gemeinsam festgelegt. jedes PM macht sich\n//    Gedanken betreffend Sinnhaftigkeit und Effizienz.\"\nconst search3 = body.search(\n  \"diskutiert und festgelegt, jedes PM macht sich Gedanken zu Sinnhaftigkeit und Effizienz.\",\n  { matchCase: true }\n);\nsearch3.load(\"items\");\nawait context.sync();\nif (search3.items.length > 0) {\n  search3.items[0].insertText(\n    \"diskutiert und gemeinsam festgelegt. jedes PM macht sich Gedanken betreffend Sinnhaftigkeit und Effizienz.\",\n    \"Replace\"\n  );\n  await context.sync();\n}\n\n// 4) Same cell, paragraph 2:\n//    \"Jedes PM erhebt Einw\u00e4nde falls es der Meinung ist, das es einen...\"\n//      -> \"Jedes PM erhebt Einw\u00e4nde, falls es der Meinung ist, dass es einen...\"\nconst search4 = body.search(\n  \"Einw\u00e4nde falls es der Meinung ist, das es einen besseren oder effizienteren\",\n  { matchCase: true }\n);\nsearch4.load(\"items\");\nawait context.sync();\nif (search4.items.length > 0) {\n  search4.items[0].insertText(\n    \"Einw\u00e4nde, falls es der Meinung ist, dass es einen besseren oder effizienteren\",\n    \"Replace\"\n  );\n  await context.sync();\n}\n\n// Note: the original document also contains Word's automatically managed\n// \"_GoBack\" bookmark (last-edit-location marker). This runtime does not\n// support deleting/relocating it through the Office.js bookmark APIs\n// (deleteBookmark silently no-ops for it), so it is intentionally left\n// untouched here rather than risk creating a duplicate bookmark name.\n", "ps1": "# Apply the four textual edits described in the commit\n# \"Anpassungen Einleitung , Organigramm\".\n#\n# Each edit is performed as a targeted Find/Replace on a unique,\n# sufficiently long substring so there is no risk of an accidental\n# match elsewhere in the document.\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text($findText, $replaceText) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    $find.Forward = $true\n    $find.Wrap = 1  # wdFindContinue\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    # wdReplaceAll = 2\n    $find.Execute($findText, $false, $true, $false, $false, $false, $true, 1, $false, $replaceText, 2) | Out-Null\n}\n\n# 1) \"Zuverl\u00e4ssigkeit\" table / Standards cell:\n#    \"Zu aufsetzende Dokumente ... abgegeben, der Zust\u00e4ndige stellt...\"\n#      -> \"Aufzusetzende Dokumente ... abgegeben. Der Zust\u00e4ndige stellt...\"\nReplace-Text `\n    \"Zu aufsetzende Dokumente werden vom jeweiligen aufgetragenen PM zum vereinbarten Zeitpunkt abgegeben, der Zust\u00e4ndige stellt zudem sicher, dass\" `\n    \"Aufzusetzende Dokumente werden vom jeweiligen aufgetragenen PM zum vereinbarten Zeitpunkt abgegeben. Der Zust\u00e4ndige stellt zudem sicher, dass\"\n\n# 2) \"Engagement\" table / Indikatoren cell:\n#    \"...fragen sie nach Hilfe und nehmen diese an.\"\n#      -> \"...fragen sie nach Hilfe und nehmen diese gegebenenfalls an.\"\nReplace-Text `\n    \"fragen sie nach Hilfe und nehmen diese an.\" `\n    \"fragen sie nach Hilfe und nehmen diese gegebenenfalls an.\"\n\n# 3) \"Engagement\" table / Standards cell, paragraph 1:\n#    \"...diskutiert und festgelegt, jedes PM macht sich Gedanken zu\n#    Sinnhaftigkeit und Effizienz.\"\n#      -> \"...diskutiert und gemeinsam festgelegt. jedes PM macht sich\n#    Gedanken betreffend Sinnhaftigkeit und Effizienz.\"\nReplace-Text `\n    \"diskutiert und festgelegt, jedes PM macht sich Gedanken zu Sinnhaftigkeit und Effizienz.\" `\n    \"diskutiert und gemeinsam festgelegt. jedes PM macht sich Gedanken betreffend Sinnhaftigkeit und Effizienz.\"\n\n# 4) Same cell, paragraph 2:\n#    \"Jedes PM erhebt Einw\u00e4nde falls es der Meinung ist, das es einen...\"\n#      -> \"Jedes PM erhebt Einw\u00e4nde, falls es der Meinung ist, dass es einen...\"\nReplace-Text `\n    \"Einw\u00e4nde falls es der Meinung ist, das es einen besseren oder effizienteren\" `\n    \"Einw\u00e4nde, falls es der Meinung ist, dass es einen besseren oder effizienteren\"\n"}
